$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9117641448974609
$ws.Range("B1").Value = 1.465642094612122
$ws.Range("D1").Value = 1.657983422279358
$ws.Range("E1").Value = 1.083024621009827
